$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old"/"_new" header suffixes to the input-file-specific suffixes
# (columns A:J keep the "FV2210" formatversion, columns L:U the "FV2304" one;
# column K ("diff") is left untouched).
$fv2210Headers = @(
    "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210",
    "Segment ID_FV2210", "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210", "Bedingung_FV2210"
)
$fv2304Headers = @(
    "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304",
    "Segment ID_FV2304", "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a proper Excel table.
$tableRange = $ws.Range("A1:U63")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""
